# "Generate Report for Handoff"
#
# The localization status report previously reflected a handback state
# ("Handed back: in sync with en-US") together with the timestamp of that
# handback. Re-running the report for a fresh handoff updates the status
# text to "Ready for handoff" and refreshes the associated timestamps on
# every sheet that surfaces them (the Overview roll-up plus each
# per-language detail sheet). The Status column (and its mirrored summary
# columns on the Overview sheet) is then re-sized to fit the new, shorter
# status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhCn     = $wb.Worksheets.Item("zh-cn")
$deDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-15 14:56:45"

# --- zh-cn detail sheet ----------------------------------------------
$zhCn.Range("C2").Value = $newStatus
$zhCn.Range("H2").Value = "2016-08-15 14:56:40"

# --- de-de detail sheet ----------------------------------------------
$deDe.Range("C2").Value = $newStatus
$deDe.Range("H2").Value = "2016-08-15 14:56:45"

# --- Resize the Status columns (and their Overview mirrors) now that
#     the text is shorter. Columns are addressed by their 1-based index
#     (E=5, F=6, C=3) rather than by letter. ---------------------------
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332
$zhCn.Columns.Item(3).ColumnWidth     = 16.333333333333332
$deDe.Columns.Item(3).ColumnWidth     = 16.333333333333332
